$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to remain text so numeric-looking values
# ("24.10", "0.0000118", etc.) are not silently coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.035.26"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "3.566.95"
$ws.Range("E3").Value = "  +4.73%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "600.64"
$ws.Range("E5").Value = "  +3.42%  "

$ws.Range("D6").Value = "138.19"
$ws.Range("E6").Value = "  +4.10%  "

$ws.Range("D7").Value = "3.565.75"
$ws.Range("E7").Value = "  +4.68%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  +3.61%  "

$ws.Range("E10").Value = "  +3.46%  "

$ws.Range("D11").Value = "6.98"
$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("E12").Value = "  +4.44%  "

$ws.Range("D13").Value = "4.174.39"
$ws.Range("E13").Value = "  +4.76%  "

$ws.Range("E14").Value = "  +4.03%  "

$ws.Range("D15").Value = "27.35"
$ws.Range("E15").Value = "  +5.46%  "

$ws.Range("D16").Value = "3.570.10"
$ws.Range("E16").Value = "  +5.11%  "

$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").Value = "65.022.78"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").Value = "10.13"
$ws.Range("E19").Value = "  +8.65%  "

$ws.Range("D20").Value = "14.44"
$ws.Range("E20").Value = "  +8.00%  "

$ws.Range("D21").Value = "5.88"
$ws.Range("E21").Value = "  +4.09%  "

$ws.Range("D22").Value = "390.86"
$ws.Range("E22").Value = "  +3.25%  "

$ws.Range("D23").Value = "0.579"
$ws.Range("E23").Value = "  +7.70%  "

$ws.Range("D24").Value = "3.712.58"
$ws.Range("E24").Value = "  +4.76%  "

$ws.Range("D25").Value = "74.16"
$ws.Range("E25").Value = "  +3.59%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").Value = "0.0000118"
$ws.Range("E27").Value = "  +13.88%  "

$ws.Range("D28").Value = "7.72"
$ws.Range("E28").Value = "  +8.21%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  +5.86%  "

$ws.Range("D31").Value = "8.36"
$ws.Range("E31").Value = "  +5.96%  "

$ws.Range("D32").Value = "3.576.67"
$ws.Range("E32").Value = "  +4.38%  "

$ws.Range("D33").Value = "1.44"
$ws.Range("E33").Value = "  +22.48%  "

$ws.Range("D34").Value = "24.10"
$ws.Range("E34").Value = "  +5.93%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = "0.146"
$ws.Range("E36").Value = "  +2.72%  "

$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("D38").Value = "6.98"
$ws.Range("E38").Value = "  +6.28%  "

$ws.Range("E39").Value = "  +8.43%  "

$ws.Range("E40").Value = "  +10.35%  "

$ws.Range("E41").Value = "  +7.63%  "

$ws.Range("E42").Value = "  +3.96%  "

$ws.Range("D43").Value = "26.99"
$ws.Range("E43").Value = "  +22.17%  "

$ws.Range("D44").Value = "42.68"
$ws.Range("E44").Value = "  +1.94%  "

$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").Value = "4.48"
$ws.Range("E46").Value = "  +6.01%  "

$ws.Range("E47").Value = "  +11.35%  "

$ws.Range("E48").Value = "  +4.47%  "

$ws.Range("D49").Value = "2.474.67"
$ws.Range("E49").Value = "  +12.80%  "

$ws.Range("D50").Value = "6.92"
$ws.Range("E50").Value = "  +7.38%  "

$ws.Range("E51").Value = "  +17.48%  "

# Restore the original (default) cell style now that the text is committed,
# so no stray number-format styling is left behind on the Price column.
$ws.Range("D2:D51").Style = "Normal"
